$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-21 down to 3-22),
# matching the weekly roll-forward of the IPO tracking table.
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new entry (이엔셀)
$ws.Range("A2").Value = "이엔셀"
$ws.Range("B2").Value = "2024.06.17~06.21"
$ws.Range("C2").Value = "13,600~15,300"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "21308"
$ws.Range("F2").Value = "NH투자증권"

# The oldest entry (노브랜드, previously row 21, now pushed to row 22) drops
# off the bottom of the table.
$ws.Rows("22:22").Delete()
